$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the old "Total" row (currently row 22) down to row 24, leaving
# rows 21-23 empty for the two new data rows plus one blank spacer row
# (matching rows 9, 12, 15 elsewhere in the sheet).
$ws.Rows(22).Insert()
$ws.Rows(23).Insert()

# Shared strings are appended in first-reference order, so touch A22
# ("Collections Home Page") before A21 ("Sorted graphics by collection")
# to land them at shared-string indices 20 and 21 respectively.

# New row 22: "Collections Home Page"
$ws.Cells.Item(20, 2).Copy($ws.Cells.Item(22, 2))
$ws.Cells.Item(20, 3).Copy($ws.Cells.Item(22, 3))
$ws.Cells.Item(22, 1).Value = "Collections Home Page"
$ws.Cells.Item(22, 2).Value = 43368
$ws.Cells.Item(22, 3).Value = 2.5

# New row 21: "Sorted graphics by collection"
$ws.Cells.Item(20, 2).Copy($ws.Cells.Item(21, 2))
$ws.Cells.Item(20, 3).Copy($ws.Cells.Item(21, 3))
$ws.Cells.Item(21, 1).Value = "Sorted graphics by collection"
$ws.Cells.Item(21, 2).Value = 43367
$ws.Cells.Item(21, 3).Value = 1

# Extend the Total formula (now on row 24) to include the new rows.
$ws.Range("C24").Formula = "=SUM(C2:C23)"

# Update the selection to match the new last-used cell.
[void]$ws.Range("C23").Select()
